$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Fix role spelling/spacing: "principalInvestigator" -> "principal Investigator"
$ws.Range("G5").Value2 = "principal Investigator"

# Fix role spelling/spacing: "metadataProvider" -> "metadata Provider"
$ws.Range("G7").Value2 = "metadata Provider"
$ws.Range("G8").Value2 = "metadata Provider"

# Make Personnel sheet the active sheet/tab, with F18 selected
$ws.Activate() | Out-Null
$ws.Range("F18").Select() | Out-Null
